# Updated symbol list on Fri Dec 23 08:40:12 UTC 2022 with GitHub Actions
# Refresh crypto price table: update Price column figures and re-sort a block
# of exchange-token rows ("One" moves up, the rest shift down one row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.13"
$ws.Range("D3").Value = "'21.94"
$ws.Range("D4").Value = "'5.376"
$ws.Range("D5").Value = "'0.05800"
$ws.Range("D6").Value = "'3.368"
$ws.Range("D7").Value = "'6.328"
$ws.Range("D8").Value = "'0.8084"
$ws.Range("D9").Value = "'1.006"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.0005892"
$ws.Range("E10").Value = "9OneONE"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1427"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07505"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03186"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03027"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'4.160"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D16").Value = "'0.09386"
$ws.Range("E16").Value = "15BitMartTokenBMX"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001588"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04821"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").Value = "'0.005644"
$ws.Range("D20").Value = "'0.004092"
$ws.Range("D21").Value = "'0.0009940"
$ws.Range("D23").Value = "'3.703"
$ws.Range("D24").Value = "'2.245"
$ws.Range("D25").Value = "'0.3256"
$ws.Range("D27").Value = "'0.0003999"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"
$ws.Range("D40").Value = "'0.03881"
$ws.Range("D41").Value = "'0.006352"
$ws.Range("D42").Value = "'0.1072"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.006688"
$ws.Range("D48").Value = "'0.1447"
